$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values, forcing text format to preserve exact
# string representation (e.g. "2.500", "1.003") without Excel converting
# them to numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.721.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2783"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06864"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07559"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6282"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.049.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009312"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.698.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.491"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.876"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.861"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1276"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.445"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06207"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.421"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.794"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.763"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.725"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.059"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6425"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.500"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.724"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01709"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.430"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.142.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8678"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.966.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.595"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.401"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05472"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4494"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) column (E) values; these are never valid numbers
# (padded with spaces), so plain text assignment is safe.
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  -7.73%  "
$ws.Range("E17").Value = "  -4.66%  "
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  -6.77%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -7.06%  "
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  -2.31%  "
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -6.16%  "
$ws.Range("E42").Value = "  -6.79%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  -6.05%  "
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  -1.53%  "
